$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": a brand-new order (#6, Sagar Borse, 18:47) came in.
#     It goes at the top of the order list, right under the header, pushing
#     every existing order down by one row. ---
$ws = $wb.Worksheets.Item("All Orders")

# Make room for the new order: shift rows 2..6 down to 3..7
$ws.Rows("2:2").Insert()

# Populate the new order in row 2
$ws.Cells.Item(2, 1).Value = 6
$ws.Cells.Item(2, 2).Value = "2026-01-13 18:47"
$ws.Cells.Item(2, 3).Value = "Sagar Borse"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "7588930329"
$ws.Cells.Item(2, 5).Value = "Test,"
$ws.Cells.Item(2, 6).Value = "Girl Haldi Kunku Set x1"
$ws.Cells.Item(2, 7).Value = 25
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""

# --- Sheet "Daily Summary": roll the new order into 2026-01-13's totals ---
$ds = $wb.Worksheets.Item("Daily Summary")
$ds.Cells.Item(2, 2).Value = 6    # Total Orders: 5 -> 6
$ds.Cells.Item(2, 5).Value = 25   # Revenue: 0 -> 25
$ds.Cells.Item(2, 7).Value = 25   # Pending: 0 -> 25
